$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, shifting existing rows 86..124 down to 87..125
$ws.Rows.Item(86).Insert()

# Populate the new row 86 with the new data record
$ws.Cells.Item(86, 1).Value = 5
$ws.Cells.Item(86, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(86, 3).Value = "Maule"
$ws.Cells.Item(86, 4).Value = 45205
$ws.Cells.Item(86, 5).Value = 7
$ws.Cells.Item(86, 6).Value = 100112026
$ws.Cells.Item(86, 7).Value = "Haba"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 400
$ws.Cells.Item(86, 11).Value = 9000
$ws.Cells.Item(86, 12).Value = 9000
$ws.Cells.Item(86, 13).Value = 9000
$ws.Cells.Item(86, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(86, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(86, 16).Value = 360
$ws.Cells.Item(86, 17).Value = 25
$ws.Cells.Item(86, 18).Value = "Hortaliza"
